$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff. Column D sometimes holds numeric-looking text
# (e.g. "74.00", "9.00") that must stay TEXT (trailing zeros matter).
# Excel auto-converts a numeric-looking Value to a real number, so for
# those we prefix with a literal single-quote (forces text entry, same
# as typing into Excel) and then ClearFormats() to drop the transient
# quote-prefix style Excel stamps on the cell, keeping the original
# (unstyled) appearance intact.

$ws.Range("D2").Value = "41.743.91"
$ws.Range("E2").Value = "  -4.71%  "
$ws.Range("D3").Value = "2.212.66"
$ws.Range("E3").Value = "  -5.72%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'246.19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("D6").Value = "'0.628"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.98%  "
$ws.Range("D7").Value = "'69.99"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.65%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.549"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -7.45%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0952"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "'36.56"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.82%  "
$ws.Range("D12").Value = "'58.04"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'6.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.02%  "
$ws.Range("D15").Value = "2.542.43"
$ws.Range("E15").Value = "  -5.76%  "
$ws.Range("D16").Value = "'14.83"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.14%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -6.56%  "
$ws.Range("D18").Value = "2.213.80"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("D19").Value = "41.646.91"
$ws.Range("E19").Value = "  -4.77%  "
$ws.Range("D20").Value = "0.0₃0956"
$ws.Range("E20").Value = "  -6.32%  "
$ws.Range("D21").Value = "'74.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.89%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.50%  "
$ws.Range("D23").Value = "'235.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.03%  "
$ws.Range("E24").Value = "  +12.24%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.46"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("E29").Value = "  -4.95%  "
$ws.Range("D30").Value = "'170.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "'20.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.82%  "
$ws.Range("D32").Value = "'0.119"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("D34").Value = "'0.0717"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").Value = "'5.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("E36").Value = "  -8.17%  "
$ws.Range("D37").Value = "'3.86"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").Value = "'23.21"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +19.11%  "
$ws.Range("E39").Value = "  -4.91%  "
$ws.Range("D40").Value = "'0.0274"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  -8.36%  "
$ws.Range("D42").Value = "'65.50"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'9.00"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.71%  "
$ws.Range("D45").Value = "'0.192"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("E46").Value = "  -5.25%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'4.56"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.71%  "
$ws.Range("D49").Value = "'10.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +8.50%  "
$ws.Range("E50").Value = "  -3.44%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.95%  "
